# Weighted Evaluation Table for Project Selection.xlsx
#
# The original sheet has three stacked "evaluator" blocks (header row +
# 9 criteria rows, the last one "Multidisciplinary"), living at rows
# 3-12, 14-23 and 25-34. This edit adds a brand-new criterion, "Fun",
# as an extra row placed right before the existing "Multidisciplinary"
# row in each of the three blocks - i.e. at (pre-shift) rows 12, 23
# and 34.
#
# We insert the rows from the bottom up so that the row numbers for the
# not-yet-processed blocks don't need to be recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 3 (rows 25-34 originally) -> insert new row before row 34 ---
$ws.Rows("34:34").Insert()
$ws.Range("A34").Value = "Fun"
$ws.Range("G34").Value = "Fun"

# --- Block 2 (rows 14-23 originally) -> insert new row before row 23 ---
$ws.Rows("23:23").Insert()
$ws.Range("A23").Value = "Fun"
$ws.Range("G23").Value = "Fun"

# --- Block 1 (rows 3-12 originally) -> insert new row before row 12 ---
$ws.Rows("12:12").Insert()
$ws.Range("A12").Value = "Fun"
$ws.Range("G12").Value = "Fun"

# Restore the view: scrolled down to the bottom block, with G36 selected
# (mirrors the author's on-screen state captured in the saved file).
[void]$ws.Range("G36").Select()
